$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated TPM-derived values for the Fgf1-Cspg4 LR-pair sheet

# Row 2
$ws.Range("G2").Value = 1.729797666666667
$ws.Range("H2").Value = 5.189393000000001
$ws.Range("I2").Value = 0.06436583050179444
$ws.Range("J2").Value = 0.06436583050179444
$ws.Range("M2").Value = 2.423077
$ws.Range("N2").Value = 7.269231
$ws.Range("O2").Value = 0.03932631260408408
$ws.Range("P2").Value = 0.03932631260408408
$ws.Range("Q2").Value = 4.191432940753667
$ws.Range("R2").Value = 37.72289646678301
$ws.Range("S2").Value = 0.002531270771335058
$ws.Range("T2").Value = 0.002531270771335058

# Row 3
$ws.Range("G3").Value = 1.729797666666667
$ws.Range("H3").Value = 5.189393000000001
$ws.Range("I3").Value = 0.06436583050179444
$ws.Range("J3").Value = 0.06436583050179444
$ws.Range("O3").Value = 0.2611559628478186
$ws.Range("P3").Value = 0.2611559628478186
$ws.Range("Q3").Value = 27.83423191425556
$ws.Range("R3").Value = 250.5080872283
$ws.Range("S3").Value = 0.01680952043919562
$ws.Range("T3").Value = 0.01680952043919562

# Row 4
$ws.Range("G4").Value = 1.729797666666667
$ws.Range("H4").Value = 5.189393000000001
$ws.Range("I4").Value = 0.06436583050179444
$ws.Range("J4").Value = 0.06436583050179444
$ws.Range("O4").Value = 0.6995177245480974
$ws.Range("P4").Value = 0.6995177245480974
$ws.Range("Q4").Value = 74.55521352407335
$ws.Range("R4").Value = 670.9969217166602
$ws.Range("S4").Value = 0.04502503929126377
$ws.Range("T4").Value = 0.04502503929126377

# Row 5
$ws.Range("I5").Value = 0.2200595722726403
$ws.Range("J5").Value = 0.2200595722726403
$ws.Range("M5").Value = 2.423077
$ws.Range("N5").Value = 7.269231
$ws.Range("O5").Value = 0.03932631260408408
$ws.Range("P5").Value = 0.03932631260408408
$ws.Range("Q5").Value = 14.33004022415266
$ws.Range("R5").Value = 128.970362017374
$ws.Range("S5").Value = 0.008654131530714886
$ws.Range("T5").Value = 0.008654131530714886

# Row 6
$ws.Range("I6").Value = 0.2200595722726403
$ws.Range("J6").Value = 0.2200595722726403
$ws.Range("O6").Value = 0.2611559628478186
$ws.Range("P6").Value = 0.2611559628478186
$ws.Range("S6").Value = 0.0574698694807405
$ws.Range("T6").Value = 0.0574698694807405

# Row 7
$ws.Range("I7").Value = 0.2200595722726403
$ws.Range("J7").Value = 0.2200595722726403
$ws.Range("O7").Value = 0.6995177245480974
$ws.Range("P7").Value = 0.6995177245480974
$ws.Range("S7").Value = 0.1539355712611849
$ws.Range("T7").Value = 0.1539355712611849

# Row 8
$ws.Range("H8").Value = 57.69206699999999
$ws.Range("I8").Value = 0.7155745972255653
$ws.Range("J8").Value = 0.7155745972255653
$ws.Range("M8").Value = 2.423077
$ws.Range("N8").Value = 7.269231
$ws.Range("O8").Value = 0.03932631260408408
$ws.Range("P8").Value = 0.03932631260408408
$ws.Range("Q8").Value = 46.59744021005299
$ws.Range("R8").Value = 419.3769618904769
$ws.Range("S8").Value = 0.02814091030203414
$ws.Range("T8").Value = 0.02814091030203414

# Row 9
$ws.Range("H9").Value = 57.69206699999999
$ws.Range("I9").Value = 0.7155745972255653
$ws.Range("J9").Value = 0.7155745972255653
$ws.Range("O9").Value = 0.2611559628478186
$ws.Range("P9").Value = 0.2611559628478186
$ws.Range("Q9").Value = 309.4416577219666
$ws.Range("S9").Value = 0.1868765729278825
$ws.Range("T9").Value = 0.1868765729278825

# Row 10
$ws.Range("H10").Value = 57.69206699999999
$ws.Range("I10").Value = 0.7155745972255653
$ws.Range("J10").Value = 0.7155745972255653
$ws.Range("O10").Value = 0.6995177245480974
$ws.Range("P10").Value = 0.6995177245480974
$ws.Range("Q10").Value = 828.85308047206
$ws.Range("R10").Value = 7459.67772424854
$ws.Range("S10").Value = 0.5005571139956487
$ws.Range("T10").Value = 0.5005571139956487

Write-Host "Updated 82 cells"
